$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly cryptos-list refresh: update Price (D) and Volume(1h) (E) columns
# plus the Kaspa/EthereumClassic row swap (rows 36-37).

$ws.Range("D2").Value = '76.150.79'
$ws.Range("E2").Value = '  -1.21%  '
$ws.Range("D3").Value = '3.065.78'
$ws.Range("E3").Value = '  +3.53%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '''197.88'
$ws.Range("E5").Value = '  -1.16%  '
$ws.Range("D6").Value = '''615.21'
$ws.Range("E6").Value = '  +3.00%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.40%  '
$ws.Range("E9").Value = '  +5.77%  '
$ws.Range("D10").Value = '3.063.09'
$ws.Range("E10").Value = '  +3.58%  '
$ws.Range("D11").Value = '''0.439'
$ws.Range("E11").Value = '  -0.74%  '
$ws.Range("E12").Value = '  -0.18%  '
$ws.Range("D13").Value = '''5.21'
$ws.Range("E13").Value = '  +6.55%  '
$ws.Range("E14").Value = '  +3.04%  '
$ws.Range("D15").Value = '''28.89'
$ws.Range("E15").Value = '  +1.62%  '
$ws.Range("D16").Value = '76.224.42'
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("E17").Value = '  +2.14%  '
$ws.Range("D18").Value = '3.063.51'
$ws.Range("E18").Value = '  +3.65%  '
$ws.Range("D19").Value = '''13.55'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").Value = '''9.10'
$ws.Range("E20").Value = '  +4.53%  '
$ws.Range("D21").Value = '''379.22'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("E22").Value = '  +8.69%  '
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").Value = '3.231.77'
$ws.Range("E24").Value = '  +3.71%  '
$ws.Range("D25").Value = '''72.03'
$ws.Range("E25").Value = '  -1.21%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  -0.14%  '
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("E32").Value = '  +2.03%  '
$ws.Range("D33").Value = '''497.53'
$ws.Range("E33").Value = '  -0.29%  '
$ws.Range("E34").Value = '  +4.15%  '
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '''20.67'
$ws.Range("E36").Value = '  +2.26%  '
$ws.Range("B37").Value = 'Kaspa'
$ws.Range("C37").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D37").Value = '''0.123'
$ws.Range("E37").Value = '  +10.53%  '
$ws.Range("D38").Value = '''162.90'
$ws.Range("E38").Value = '  -2.07%  '
$ws.Range("D39").Value = '''20.04'
$ws.Range("E39").Value = '  +0.47%  '
$ws.Range("D40").Value = '''193.19'
$ws.Range("E40").Value = '  +7.24%  '
$ws.Range("E41").Value = '  -5.81%  '
$ws.Range("E42").Value = '  -9.08%  '
$ws.Range("D44").Value = '''0.792'
$ws.Range("E44").Value = '  +20.21%  '
$ws.Range("D45").Value = '''5.09'
$ws.Range("E45").Value = '  +2.96%  '
$ws.Range("E46").Value = '  +4.00%  '
$ws.Range("E47").Value = '  +2.69%  '
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("E49").Value = '  +4.74%  '
$ws.Range("D50").Value = '''0.593'
$ws.Range("E50").Value = '  +0.00%  '
$ws.Range("E51").Value = '  -1.07%  '

$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D50").ClearFormats()
